$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D; this shifts existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy number formats/styles from column E (the old column D, now shifted) into
# the new column D, for each contiguous block of data rows, so the new column
# matches the row's style (date format for header rows, number format for data rows).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the new column D with the latest period's financial figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 530900
$ws.Range("D9").Value = 226900
$ws.Range("D10").Value = 304000
$ws.Range("D12").Value = 61500
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 37200
$ws.Range("D15").Value = 22600
$ws.Range("D17").Value = 555500
$ws.Range("D18").Value = -24600
$ws.Range("D20").Value = -900
$ws.Range("D21").Value = 8400
$ws.Range("D22").Value = 6800
$ws.Range("D23").Value = -32300
$ws.Range("D24").Value = -17900
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -14300
$ws.Range("D27").Value = -14300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -8600
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 900
$ws.Range("D33").Value = -22900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -22900
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 56400
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 127000
$ws.Range("D44").Value = 79700
$ws.Range("D45").Value = 22600
$ws.Range("D46").Value = 285800
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 22900
$ws.Range("D49").Value = 287100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 42400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 638100
$ws.Range("D57").Value = 28800
$ws.Range("D58").Value = 35000
$ws.Range("D59").Value = 69600
$ws.Range("D60").Value = 133400
$ws.Range("D61").Value = 69500
$ws.Range("D62").Value = 36800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 239700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 102300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 398400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -22900
$ws.Range("D83").Value = 33900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 33000
$ws.Range("D91").Value = -7900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -8400
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -49500
$ws.Range("D101").Value = -7700
$ws.Range("D102").Value = -32600

# Row 91 (Capital Expenditures) received corrected historical figures in
# addition to the new column, so update the shifted columns F:J explicitly.
$ws.Range("F91").Value = -3200
$ws.Range("G91").Value = -4100
$ws.Range("H91").Value = -4200
$ws.Range("I91").Value = -1800
$ws.Range("J91").Value = -2200
